# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" tracking sheet
# and moves the special "latest row" date formatting from the old last row
# (row 89) onto the new last row (row 90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (89) loses the distinct "last row" date format and
# reverts to the regular date/time format used by all the other data rows.
$ws.Cells.Item(89, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 90.
$ws.Cells.Item(90, 1).Value = 45677
$ws.Cells.Item(90, 2).Value = 212
$ws.Cells.Item(90, 3).Value = 209
$ws.Cells.Item(90, 4).Value = 211

# The new last row takes on the distinct "last row" date-only format.
$ws.Cells.Item(90, 1).NumberFormat = "YYYY-MM-DD"
